# Penalty/Reward System (unfinished) - data refresh for forecast_summary_B0CNR6L1DH
#
# Shifts the weekly forecast window forward by one week on the
# "Forecast Comparison" sheet (Week_Start_Date + MyForecast columns) and
# refreshes the derived statistics on the "Summary" sheet to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value to a cell as plain text, even when the text
# looks like a date/number, without leaving a lingering custom style
# behind (NumberFormat "@" forces text entry; resetting the Style back
# to "Normal" afterwards drops the cell back to the default style so
# only the value itself changes).
# ---------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet 1: "Forecast Comparison"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

Set-TextValue $ws1.Range("B2")  "2025-01-12"
$ws1.Range("D2").Value = 20

Set-TextValue $ws1.Range("B3")  "2025-01-19"
$ws1.Range("D3").Value = 21

Set-TextValue $ws1.Range("B4")  "2025-01-26"
$ws1.Range("D4").Value = 21

Set-TextValue $ws1.Range("B5")  "2025-02-02"
$ws1.Range("D5").Value = 21

Set-TextValue $ws1.Range("B6")  "2025-02-09"
$ws1.Range("D6").Value = 24

Set-TextValue $ws1.Range("B7")  "2025-02-16"
$ws1.Range("D7").Value = 27

Set-TextValue $ws1.Range("B8")  "2025-02-23"
$ws1.Range("D8").Value = 20

Set-TextValue $ws1.Range("B9")  "2025-03-02"
$ws1.Range("D9").Value = 20

Set-TextValue $ws1.Range("B10") "2025-03-09"
$ws1.Range("D10").Value = 28

Set-TextValue $ws1.Range("B11") "2025-03-16"
$ws1.Range("D11").Value = 28

Set-TextValue $ws1.Range("B12") "2025-03-23"
$ws1.Range("D12").Value = 26

Set-TextValue $ws1.Range("B13") "2025-03-30"
$ws1.Range("D13").Value = 23

Set-TextValue $ws1.Range("B14") "2025-04-06"
# D14 (MyForecast) is unchanged - stays 20

Set-TextValue $ws1.Range("B15") "2025-04-13"
$ws1.Range("D15").Value = 20

Set-TextValue $ws1.Range("B16") "2025-04-20"
$ws1.Range("D16").Value = 24

Set-TextValue $ws1.Range("B17") "2025-04-27"
$ws1.Range("D17").Value = 18

# ---------------------------------------------------------------------
# Sheet 2: "Summary"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

Set-TextValue $ws2.Range("B2")  "2023-12-10 to 2025-01-05"   # Historical Range
Set-TextValue $ws2.Range("B4")  "36"                          # Max Sales
Set-TextValue $ws2.Range("B6")  "16"                          # Median Sales
Set-TextValue $ws2.Range("B8")  "860 units"                   # Total Historical Sales
Set-TextValue $ws2.Range("B9")  "361"                         # Total Forecast (16 Weeks)
Set-TextValue $ws2.Range("B10") "173"                         # Total Forecast (8 Weeks)
Set-TextValue $ws2.Range("B11") "83"                          # Total Forecast (4 Weeks)
Set-TextValue $ws2.Range("B12") "28"                          # Max Forecast
Set-TextValue $ws2.Range("B13") "2025-03-09"                  # Max Forecast Week
Set-TextValue $ws2.Range("B14") "18"                          # Min Forecast
Set-TextValue $ws2.Range("B15") "2025-04-27"                  # Min Forecast Week
